# Auto-generated Excel COM-interop script to apply scheduled-runner market-data refresh
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H..N) across all 8 sheets.
$wb = $excel.ActiveWorkbook

# ==== Sheet: ALC ====
$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 73.59999999999999
$ws.Range("I11").Value = 73.59999999999999
$ws.Range("K11").Value = 73.59999999999999
$ws.Range("M11").Value = 66.40000000000001

# Row 62
$ws.Range("H62").Value = 18520984
$ws.Range("I62").Value = 37038304
$ws.Range("J62").Value = 3663.3333
$ws.Range("K62").Value = 37038304
$ws.Range("L62").Value = 3663.3333
$ws.Range("M62").Value = -37037680
$ws.Range("N62").Value = -4911.3333

# Row 65
$ws.Range("H65").Value = 18520984
$ws.Range("I65").Value = 37038304
$ws.Range("J65").Value = 3663.3333
$ws.Range("K65").Value = 185191520
$ws.Range("L65").Value = 18316.6665
$ws.Range("M65").Value = -185188400
$ws.Range("N65").Value = -24556.6665

# Row 82
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

# Row 85
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

# Row 114
$ws.Range("H114").Value = 38000
$ws.Range("J114").Value = 38000
$ws.Range("L114").Value = 38000
$ws.Range("N114").Value = -46678

# Row 132
$ws.Range("H132").Value = 8138301
$ws.Range("I132").Value = 10106468
$ws.Range("J132").Value = 19614
$ws.Range("K132").Value = 30319404
$ws.Range("L132").Value = 58842
$ws.Range("M132").Value = -30316874
$ws.Range("N132").Value = -63902

# Row 137
$ws.Range("H137").Value = 1167.9254
$ws.Range("I137").Value = 922.13635
$ws.Range("J137").Value = 1638.1305
$ws.Range("K137").Value = 2766.40905
$ws.Range("L137").Value = 4914.3915
$ws.Range("M137").Value = -216.4090500000002
$ws.Range("N137").Value = -10014.3915

# Row 138
$ws.Range("H138").Value = 1522.13
$ws.Range("I138").Value = 776.7857
$ws.Range("J138").Value = 1811.9861
$ws.Range("K138").Value = 2330.3571
$ws.Range("L138").Value = 5435.9583
$ws.Range("M138").Value = 2809.6429
$ws.Range("N138").Value = -15715.9583

# ==== Sheet: ARM ====
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 771.1667
$ws.Range("I2").Value = 497.03705
$ws.Range("J2").Value = 1593.5555
$ws.Range("K2").Value = 497.03705
$ws.Range("L2").Value = 1593.5555
$ws.Range("M2").Value = -384.03705
$ws.Range("N2").Value = -1819.5555

# Row 32
$ws.Range("H32").Value = 7686.1885
$ws.Range("I32").Value = 6400.3716
$ws.Range("J32").Value = 13686.667
$ws.Range("K32").Value = 6400.3716
$ws.Range("L32").Value = 13686.667
$ws.Range("M32").Value = -6113.3716
$ws.Range("N32").Value = -14260.667

# Row 80
$ws.Range("H80").Value = 36750
$ws.Range("J80").Value = 36750
$ws.Range("L80").Value = 36750
$ws.Range("N80").Value = -38746

# Row 83
$ws.Range("H83").Value = 36750
$ws.Range("J83").Value = 36750
$ws.Range("L83").Value = 110250
$ws.Range("N83").Value = -120234

# Row 116
$ws.Range("H116").Value = 771.1667
$ws.Range("I116").Value = 497.03705
$ws.Range("J116").Value = 1593.5555
$ws.Range("K116").Value = 497.03705
$ws.Range("L116").Value = 1593.5555
$ws.Range("M116").Value = 1796.96295
$ws.Range("N116").Value = -6181.5555

# Row 121
$ws.Range("H121").Value = 37750
$ws.Range("J121").Value = 37750
$ws.Range("L121").Value = 37750
$ws.Range("N121").Value = -41244

# Row 122
$ws.Range("H122").Value = 2036.5
$ws.Range("I122").Value = 1817.6428
$ws.Range("J122").Value = 2419.5
$ws.Range("K122").Value = 5452.928400000001
$ws.Range("L122").Value = 7258.5
$ws.Range("M122").Value = -3002.928400000001
$ws.Range("N122").Value = -12158.5

# ==== Sheet: BSM ====
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 771.1667
$ws.Range("I3").Value = 497.03705
$ws.Range("J3").Value = 1593.5555
$ws.Range("K3").Value = 497.03705
$ws.Range("L3").Value = 1593.5555
$ws.Range("M3").Value = -383.03705
$ws.Range("N3").Value = -1821.5555

# Row 20
$ws.Range("H20").Value = 3842.8572
$ws.Range("I20").Value = 3733.3333
$ws.Range("J20").Value = 4500
$ws.Range("K20").Value = 3733.3333
$ws.Range("L20").Value = 4500
$ws.Range("M20").Value = -3486.3333
$ws.Range("N20").Value = -4994

# Row 140
$ws.Range("H140").Value = 28068.889
$ws.Range("J140").Value = 28068.889
$ws.Range("L140").Value = 28068.889
$ws.Range("N140").Value = -38428.889

# ==== Sheet: CRP ====
$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 760.6957
$ws.Range("I107").Value = 416.44446
$ws.Range("K107").Value = 416.44446
$ws.Range("M107").Value = 1503.55554

# Row 122
$ws.Range("H122").Value = 714.5
$ws.Range("J122").Value = 671.3333
$ws.Range("L122").Value = 2013.9999
$ws.Range("N122").Value = -6913.9999

# Row 132
$ws.Range("H132").Value = 2097.1936
$ws.Range("I132").Value = 1736.8823
$ws.Range("J132").Value = 2534.7144
$ws.Range("K132").Value = 5210.6469
$ws.Range("L132").Value = 7604.1432
$ws.Range("M132").Value = -2680.6469
$ws.Range("N132").Value = -12664.1432

# ==== Sheet: CUL ====
$ws = $wb.Worksheets.Item("CUL")
# Row 87
$ws.Range("H87").Value = 1394.625
$ws.Range("I87").Value = 1042.8
$ws.Range("K87").Value = 3128.4
$ws.Range("M87").Value = -1880.4

# Row 90
$ws.Range("H90").Value = 1394.625
$ws.Range("I90").Value = 1042.8
$ws.Range("K90").Value = 9385.199999999999
$ws.Range("M90").Value = -3145.199999999999

# Row 138
$ws.Range("H138").Value = 2467.2703
$ws.Range("I138").Value = 2383.9333
$ws.Range("J138").Value = 2524.0908
$ws.Range("K138").Value = 7151.7999
$ws.Range("L138").Value = 7572.2724
$ws.Range("M138").Value = -2011.7999
$ws.Range("N138").Value = -17852.2724

# Row 140
$ws.Range("H140").Value = 24737.283
$ws.Range("I140").Value = 44122.375
$ws.Range("J140").Value = 3589.9092
$ws.Range("K140").Value = 132367.125
$ws.Range("L140").Value = 10769.7276
$ws.Range("M140").Value = -127187.125
$ws.Range("N140").Value = -21129.7276

# ==== Sheet: GSM ====
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 9480
$ws.Range("J80").Value = 6850
$ws.Range("L80").Value = 6850
$ws.Range("N80").Value = -8846

# Row 83
$ws.Range("H83").Value = 9480
$ws.Range("J83").Value = 6850
$ws.Range("L83").Value = 34250
$ws.Range("N83").Value = -44234

# Row 102
$ws.Range("H102").Value = 1367.8
$ws.Range("I102").Value = 1287.0322
$ws.Range("K102").Value = 1287.0322
$ws.Range("M102").Value = 334.9677999999999

# Row 113
$ws.Range("H113").Value = 1519.0834
$ws.Range("I113").Value = 1394
$ws.Range("K113").Value = 1394
$ws.Range("M113").Value = 776

# Row 121
$ws.Range("H121").Value = 50499.5
$ws.Range("J121").Value = 50499.5
$ws.Range("L121").Value = 50499.5
$ws.Range("N121").Value = -53993.5

# Row 122
$ws.Range("H122").Value = 4136.737
$ws.Range("I122").Value = 4099.875
$ws.Range("J122").Value = 4333.3335
$ws.Range("K122").Value = 12299.625
$ws.Range("L122").Value = 13000.0005
$ws.Range("M122").Value = -9849.625
$ws.Range("N122").Value = -17900.0005

# Row 126
$ws.Range("H126").Value = 2126.3157
$ws.Range("I126").Value = 1803
$ws.Range("J126").Value = 2485.5557
$ws.Range("K126").Value = 5409
$ws.Range("L126").Value = 7456.6671
$ws.Range("M126").Value = -2939
$ws.Range("N126").Value = -12396.6671

# ==== Sheet: LTW ====
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 4850
$ws.Range("J46").Value = 5660
$ws.Range("L46").Value = 5660
$ws.Range("N46").Value = -6036

# Row 61
$ws.Range("H61").Value = 922.6667
$ws.Range("I61").Value = 787.5
$ws.Range("K61").Value = 787.5
$ws.Range("M61").Value = -585.5

# Row 68
$ws.Range("H68").Value = 1276.8572
$ws.Range("I68").Value = 1283
$ws.Range("J68").Value = 1240
$ws.Range("K68").Value = 1283
$ws.Range("L68").Value = 1240
$ws.Range("M68").Value = -534
$ws.Range("N68").Value = -2738

# Row 71
$ws.Range("H71").Value = 1276.8572
$ws.Range("I71").Value = 1283
$ws.Range("J71").Value = 1240
$ws.Range("K71").Value = 6415
$ws.Range("L71").Value = 6200
$ws.Range("M71").Value = -2671
$ws.Range("N71").Value = -13688

# Row 76
$ws.Range("H76").Value = 10000
$ws.Range("J76").Value = 10000
$ws.Range("L76").Value = 10000
$ws.Range("N76").Value = -10676

# Row 79
$ws.Range("H79").Value = 10000
$ws.Range("J79").Value = 10000
$ws.Range("L79").Value = 10000
$ws.Range("N79").Value = -12340

# Row 113
$ws.Range("H113").Value = 922.6667
$ws.Range("I113").Value = 787.5
$ws.Range("K113").Value = 787.5
$ws.Range("M113").Value = 1382.5

# Row 132
$ws.Range("H132").Value = 2588.8845
$ws.Range("I132").Value = 2294.5881
$ws.Range("J132").Value = 3144.7778
$ws.Range("K132").Value = 6883.7643
$ws.Range("L132").Value = 9434.3334
$ws.Range("M132").Value = -4353.7643
$ws.Range("N132").Value = -14494.3334

# Row 133
$ws.Range("H133").Value = 45695.832
$ws.Range("J133").Value = 45695.832
$ws.Range("L133").Value = 45695.832
$ws.Range("N133").Value = -50755.832

# ==== Sheet: WVR ====
$ws = $wb.Worksheets.Item("WVR")
# Row 99
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

# Row 113
$ws.Range("H113").Value = 545.5
$ws.Range("I113").Value = 408.33334
$ws.Range("J113").Value = 751.25
$ws.Range("K113").Value = 1225.00002
$ws.Range("L113").Value = 2253.75
$ws.Range("M113").Value = 944.9999800000001
$ws.Range("N113").Value = -6593.75

# Row 122
$ws.Range("H122").Value = 13159919
$ws.Range("I122").Value = 14707956
$ws.Range("J122").Value = 1597.5
$ws.Range("K122").Value = 44123868
$ws.Range("L122").Value = 4792.5
$ws.Range("M122").Value = -44121418
$ws.Range("N122").Value = -9692.5

# Row 136
$ws.Range("H136").Value = 1450.9166
$ws.Range("I136").Value = 1230.4117
$ws.Range("K136").Value = 3691.2351
$ws.Range("M136").Value = -1141.2351

